$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.048.84"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "1.719.05"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "318.97"
$ws.Range("E5").Value = "  -2.65%  "
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "0.4671"
$ws.Range("E7").Value = "  +4.02%  "
$ws.Range("D8").Value = "0.3443"
$ws.Range("E8").Value = "  -3.61%  "
$ws.Range("D9").Value = "42.31"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "0.07311"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("D11").Value = "1.053"
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("D12").Value = "1.006"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "19.93"
$ws.Range("E13").Value = "  -4.75%  "
$ws.Range("D14").Value = "5.874"
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("D15").Value = "1.734.43"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").Value = "6.906"
$ws.Range("E16").Value = "  -4.49%  "
$ws.Range("D17").Value = "89.81"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").Value = "0.00001048"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "0.06313"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").Value = "1.008"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "16.47"
$ws.Range("E21").Value = "  -4.42%  "
$ws.Range("D22").Value = "5.625"
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("D23").Value = "27.118.54"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").Value = "10.83"
$ws.Range("E24").Value = "  -4.54%  "
$ws.Range("D25").Value = "2.128"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").Value = "157.12"
$ws.Range("E26").Value = "  -3.42%  "
$ws.Range("D27").Value = "19.49"
$ws.Range("E27").Value = "  -3.87%  "
$ws.Range("D28").Value = "1.919.21"
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("D29").Value = "2.142"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").Value = "119.27"
$ws.Range("E30").Value = "  -5.23%  "
$ws.Range("D31").Value = "1.021"
$ws.Range("E31").Value = "  -7.46%  "
$ws.Range("D32").Value = "0.09094"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").Value = "3.605"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").Value = "5.334"
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("D35").Value = "0.02205"
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").Value = "11.14"
$ws.Range("E36").Value = "  -6.27%  "
$ws.Range("D37").Value = "0.05837"
$ws.Range("E37").Value = "  -4.40%  "
$ws.Range("D38").Value = "0.1996"
$ws.Range("E38").Value = "  -4.84%  "
$ws.Range("D41").Value = "0.5950"
$ws.Range("E41").Value = "  -6.33%  "
$ws.Range("D42").Value = "1.137"
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("D43").Value = "7.510"
$ws.Range("E43").Value = "  -5.23%  "
$ws.Range("D44").Value = "3.653"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "12.59"
$ws.Range("E45").Value = "  -5.10%  "
$ws.Range("D46").Value = "0.5599"
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("D47").Value = "119.34"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").Value = "1.863"
$ws.Range("E48").Value = "  -4.94%  "
$ws.Range("D49").Value = "0.06666"
$ws.Range("E49").Value = "  -3.72%  "
$ws.Range("D50").Value = "1.089"
$ws.Range("E50").Value = "  -4.62%  "
$ws.Range("D51").Value = "1.007"
$ws.Range("E51").Value = "  +0.47%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "4.759"
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "1.406"
$ws.Range("E40").Value = "  +1.04%  "
